$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vwf"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 28.67276066666667
$ws.Cells.Item(2, 8).Value = 86.018282
$ws.Cells.Item(2, 9).Value = 0.9474462168692853
$ws.Cells.Item(2, 10).Value = 0.9474462168692853
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1176943333333333
$ws.Cells.Item(2, 14).Value = 0.353083
$ws.Cells.Item(2, 15).Value = 0.05829606481791055
$ws.Cells.Item(2, 16).Value = 0.05829606481791055
$ws.Cells.Item(2, 17).Value = 3.374621451489555
$ws.Cells.Item(2, 18).Value = 30.371593063406
$ws.Cells.Item(2, 19).Value = 0.05523238607009599
$ws.Cells.Item(2, 20).Value = 0.05523238607009599

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vwf"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 28.67276066666667
$ws.Cells.Item(3, 8).Value = 86.018282
$ws.Cells.Item(3, 9).Value = 0.9474462168692853
$ws.Cells.Item(3, 10).Value = 0.9474462168692853
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.183046666666667
$ws.Cells.Item(3, 14).Value = 3.54914
$ws.Cells.Item(3, 15).Value = 0.5859837360842607
$ws.Cells.Item(3, 16).Value = 0.5859837360842608
$ws.Cells.Item(3, 17).Value = 33.92121393083111
$ws.Cells.Item(3, 18).Value = 305.29092537748
$ws.Cells.Item(3, 19).Value = 0.5551880738999625
$ws.Cells.Item(3, 20).Value = 0.5551880738999626

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vwf"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 28.67276066666667
$ws.Cells.Item(4, 8).Value = 86.018282
$ws.Cells.Item(4, 9).Value = 0.9474462168692853
$ws.Cells.Item(4, 10).Value = 0.9474462168692853
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.718166
$ws.Cells.Item(4, 14).Value = 2.154498
$ws.Cells.Item(4, 15).Value = 0.3557201990978286
$ws.Cells.Item(4, 16).Value = 0.3557201990978286
$ws.Cells.Item(4, 17).Value = 20.59180183693733
$ws.Cells.Item(4, 18).Value = 185.326216532436
$ws.Cells.Item(4, 19).Value = 0.3370257568992267
$ws.Cells.Item(4, 20).Value = 0.3370257568992267

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Vwf"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.7569533333333333
$ws.Cells.Item(5, 8).Value = 2.27086
$ws.Cells.Item(5, 9).Value = 0.02501233070476559
$ws.Cells.Item(5, 10).Value = 0.02501233070476559
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.1176943333333333
$ws.Cells.Item(5, 14).Value = 0.353083
$ws.Cells.Item(5, 15).Value = 0.05829606481791055
$ws.Cells.Item(5, 16).Value = 0.05829606481791055
$ws.Cells.Item(5, 17).Value = 0.0890891179311111
$ws.Cells.Item(5, 18).Value = 0.8018020613799999
$ws.Cells.Item(5, 19).Value = 0.001458120452012029
$ws.Cells.Item(5, 20).Value = 0.001458120452012029

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vwf"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.7569533333333333
$ws.Cells.Item(6, 8).Value = 2.27086
$ws.Cells.Item(6, 9).Value = 0.02501233070476559
$ws.Cells.Item(6, 10).Value = 0.02501233070476559
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.183046666666667
$ws.Cells.Item(6, 14).Value = 3.54914
$ws.Cells.Item(6, 15).Value = 0.5859837360842607
$ws.Cells.Item(6, 16).Value = 0.5859837360842608
$ws.Cells.Item(6, 17).Value = 0.895511117822222
$ws.Cells.Item(6, 18).Value = 8.059600060399999
$ws.Cells.Item(6, 19).Value = 0.01465681899455361
$ws.Cells.Item(6, 20).Value = 0.01465681899455361

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vwf"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.7569533333333333
$ws.Cells.Item(7, 8).Value = 2.27086
$ws.Cells.Item(7, 9).Value = 0.02501233070476559
$ws.Cells.Item(7, 10).Value = 0.02501233070476559
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.718166
$ws.Cells.Item(7, 14).Value = 2.154498
$ws.Cells.Item(7, 15).Value = 0.3557201990978286
$ws.Cells.Item(7, 16).Value = 0.3557201990978286
$ws.Cells.Item(7, 17).Value = 0.5436181475866666
$ws.Cells.Item(7, 18).Value = 4.89256332828
$ws.Cells.Item(7, 19).Value = 0.008897391258199949
$ws.Cells.Item(7, 20).Value = 0.008897391258199949

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Vwf"
$ws.Cells.Item(8, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.8334926666666668
$ws.Cells.Item(8, 8).Value = 2.500478
$ws.Cells.Item(8, 9).Value = 0.02754145242594914
$ws.Cells.Item(8, 10).Value = 0.02754145242594913
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.1176943333333333
$ws.Cells.Item(8, 14).Value = 0.353083
$ws.Cells.Item(8, 15).Value = 0.05829606481791055
$ws.Cells.Item(8, 16).Value = 0.05829606481791055
$ws.Cells.Item(8, 17).Value = 0.09809736374155556
$ws.Cells.Item(8, 18).Value = 0.882876273674
$ws.Cells.Item(8, 19).Value = 0.001605558295802531
$ws.Cells.Item(8, 20).Value = 0.001605558295802531

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Vwf"
$ws.Cells.Item(9, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.8334926666666668
$ws.Cells.Item(9, 8).Value = 2.500478
$ws.Cells.Item(9, 9).Value = 0.02754145242594914
$ws.Cells.Item(9, 10).Value = 0.02754145242594913
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.183046666666667
$ws.Cells.Item(9, 14).Value = 3.54914
$ws.Cells.Item(9, 15).Value = 0.5859837360842607
$ws.Cells.Item(9, 16).Value = 0.5859837360842608
$ws.Cells.Item(9, 17).Value = 0.9860607209911112
$ws.Cells.Item(9, 18).Value = 8.87454648892
$ws.Cells.Item(9, 19).Value = 0.0161388431897446
$ws.Cells.Item(9, 20).Value = 0.0161388431897446

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Vwf"
$ws.Cells.Item(10, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8334926666666668
$ws.Cells.Item(10, 8).Value = 2.500478
$ws.Cells.Item(10, 9).Value = 0.02754145242594914
$ws.Cells.Item(10, 10).Value = 0.02754145242594913
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.718166
$ws.Cells.Item(10, 14).Value = 2.154498
$ws.Cells.Item(10, 15).Value = 0.3557201990978286
$ws.Cells.Item(10, 16).Value = 0.3557201990978286
$ws.Cells.Item(10, 17).Value = 0.5985860944493334
$ws.Cells.Item(10, 18).Value = 5.387274850044
$ws.Cells.Item(10, 19).Value = 0.009797050940402003
$ws.Cells.Item(10, 20).Value = 0.009797050940402003
